$wb = $excel.ActiveWorkbook

# Update the "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets, which carry duplicate data. Row 3 -> 南宁·万圣漫控嘉年华10, Row 5 -> 南宁·黑塔利亚同人ONLY
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1387
    $ws.Range("F5").Value = 72
}
